$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "batch 32" device rows appended after the existing data (rows 157-161),
# mirroring the pattern already used for every prior batch (e.g. rows 152-156
# for "batch 31"): id, name, mac_address, serial_num, (ip_address left blank),
# dspec_id, lang_code, is_active, cr_by, cr_dtimes.

$rows = @(
    @{ Row = 157; Id = 3000176; Name = "Finger Print Scanner 32"; Mac = "80-75-40-E8-CA-24"; Serial = "BS563Q2230824"; DspecId = 165 },
    @{ Row = 158; Id = 3000177; Name = "IRIS Scanner 32";         Mac = "0E-1A-14-4A-6D-3A"; Serial = "BS563Q2230825"; DspecId = 327 },
    @{ Row = 159; Id = 3000178; Name = "Web Camera 32";           Mac = "65-13-7F-0F-F7-53"; Serial = "BS563Q2230826"; DspecId = 736 },
    @{ Row = 160; Id = 3000179; Name = "Document Scanner 32";     Mac = "73-C4-DE-8E-C9-8D"; Serial = "BS563Q2230827"; DspecId = 801 },
    @{ Row = 161; Id = 3000180; Name = "Printer 32";              Mac = "EC-74-AB-E0-0F-38"; Serial = "BS563Q2230828"; DspecId = 920 }
)

# Populate column-by-column (not row-by-row) so new shared-string entries are
# interned in the same order the source workbook has them: all five names,
# then all five mac addresses, then all five serial numbers, etc.
foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Id
}
foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.Name
}
foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 3).Value = $r.Mac
}
foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 4).Value = $r.Serial
}
# Column E (ip_address) intentionally left blank, matching every other row.
foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 6).Value = $r.DspecId
}
foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 7).Value = "eng"
}
foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 8).Value = $true
    $ws.Cells.Item($r.Row, 8).HorizontalAlignment = -4131
}
foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 9).Value = "superadmin"
}
foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 10).Value = "now()"
}

# Best-effort restore of the author's on-screen view state: scrolled so row
# 113 is at the top, with the whole of column K onward selected (active cell
# sitting in row 113 of that selection).
$excel.ActiveWindow.ScrollRow = 113
$ws.Cells.Item(113, 11).Select() | Out-Null
$ws.Range("K1:XFD1048576").Select() | Out-Null
